# Add new Sheet2 after Sheet1 with the metric conversion data
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Nothing"
$ws2.Range("B1").Value = 950
$ws2.Range("C1").Value = 1000
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = 975

$ws1.Range("C1").Select() | Out-Null

$ws2.Range("A3").Select() | Out-Null

$ws2.Activate()
